# fix problems with transferring to WP
# Replace the full amenities list in column A with the curated/shortened list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "24 Hour Availability",
    "Basketball Court",
    "Boat Docks",
    "Business Center",
    "Clubhouse",
    "Coffee System",
    "Community-Wide WiFi",
    "Controlled Access",
    "Courtyard",
    "Disposal Chutes",
    "Elevator",
    "Fenced Lot",
    "Fitness Center",
    "Gameroom",
    "Garden",
    "Grill",
    "Guest Apartment",
    "Key Fob Entry",
    "Lake Access",
    "Lounge",
    "Maintenance on site",
    "Mud Room",
    "Multi Use Room",
    "On-Site ATM",
    "Online Services",
    "Package Service",
    "Pet Play Area",
    "Picnic Area",
    "Planned Social Activities",
    "Pond",
    "Pool",
    "Property Manager on Site",
    "Public Transportation",
    "Recycling",
    "Renters Insurance Program",
    "Roof Terrace",
    "Sauna",
    "Spa",
    "Storage Space",
    "Sundeck",
    "Tennis Court",
    "Trash Pickup - Curbside",
    "Volleyball Court",
    "Walking/Biking Trails",
    "Waterfront",
    "Wi-Fi at Pool and Clubhouse"
)

$oldLastRow = $ws.UsedRange.Rows.Count
$newCount = $values.Count

# Each value in the sheet is wrapped in literal single quotes, e.g. 'Pool'.
# Assigning a string like "'Pool'" directly through .Value triggers Excel's
# quote-prefix (text literal) auto-detection and strips the leading quote,
# so instead build the text with a formula using CHAR(39) for the apostrophes
# and then convert the formulas to static values via copy / paste-special.
for ($i = 0; $i -lt $newCount; $i++) {
    $cell = $ws.Cells.Item($i + 1, 1)
    $cell.Formula = "=CHAR(39)&""" + $values[$i] + """&CHAR(39)"
}

$newRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($newCount, 1))
$newRange.Copy()
$newRange.PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Remove any leftover rows beyond the new list (old sheet had 90 rows)
if ($oldLastRow -gt $newCount) {
    $deleteRange = $ws.Range($ws.Cells.Item($newCount + 1, 1), $ws.Cells.Item($oldLastRow, 1))
    $deleteRange.EntireRow.Delete()
}

$wb.Save()
